# PLANTILLA TRANSFER DIRECTO CRM.xlsx
# Reproduce the template clean-up commit:
#  - Excel principal: siempre generado desde cero -> strip the example/sample
#    data that used to live in this template (pharmacy/contact info in
#    C13:C16 and the IALOZON product lines in B21:E27) so the sheet is a
#    bare template again.
#  - Nudge/resize the HEFAME logo image slightly.
#  - Leave the cursor/selection on D5 instead of F5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Selection moves from F5 to D5 -----------------------------------
[void]$ws.Range("D5").Select()

# --- Clear the sample pharmacy data (Codigo Hefame / Nombre / Telefono /
#     Poblacion) out of the "DATOS FARMACIA" block. Styles/number formats
#     on these cells are left untouched, only the values are removed. ---
$ws.Range("C13").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("C16").ClearContents()

# --- Clear the sample IALOZON order lines (Cantidad, CN, Descripcion,
#     Descuento) from the content table, rows 21-27. ------------------
$ws.Range("B21:E27").ClearContents()

# --- Reposition/resize the HEFAME logo picture in the header. --------
$shp = $ws.Shapes.Item(1)
$shp.Left = 542.8585039370079
$shp.Top = 10.641496062992125
$shp.Width = 134.44818897637796
$shp.Height = 35.900944881889764
